$wb = $excel.ActiveWorkbook

# --- "About" sheet: trim the Source block down to "None" and refresh the Notes text ---
$ws1 = $wb.Worksheets.Item("About")

# Remove the hyperlink to the old source (sciencedirect.com link in B6)
$ws1.Hyperlinks.Delete()

# Remove rows 4-7 (year 2014, paper title, hyperlink url, "Page 190, column 2")
# which collapses what used to be rows 9-11 up to rows 5-7.
$ws1.Range("A4:A7").EntireRow.Delete()

# New note explaining the US-specific override (goes into the now-empty row 9,
# written before B3 so shared-string ordering matches)
$ws1.Range("A9").Value = "In the US, we set this to 0 so that increasing EV chargers does not induce additional deployment."

# Source reference is now "None"
$ws1.Range("B3").Value = "None"

# Remove the now-unused built-in "Hyperlink" cell style
$wb.Styles.Item("Hyperlink").Delete()

# --- "EoCSoEVMS" sheet: zero out the US market-share-change coefficient ---
$ws2 = $wb.Worksheets.Item("EoCSoEVMS")
$ws2.Range("B2").Value = 0

# --- restore selections so "About" ends up the active tab ---
$ws2.Activate()
$ws2.Range("B3").Select()

$ws1.Activate()
$ws1.Range("A4:XFD7").Select()
